# Generate Report for Archive
#
# 1. Update the localization status text from "Ready for handoff" to
#    "In Translation" everywhere it appears (Overview!E2/F2, zh-cn!C2,
#    de-de!C2 all share the same string).
# 2. Narrow the "Status" column width on the Overview sheet (columns E/F)
#    and on the zh-cn / de-de sheets (column C) to match the new report
#    layout.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Update status strings -------------------------------------------------
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# --- Update column widths ---------------------------------------------------
$overview.Range("E1").ColumnWidth = 13.4101848602295
$overview.Range("F1").ColumnWidth = 13.4101848602295

$zhcn.Range("C1").ColumnWidth = 13.4101848602295
$dede.Range("C1").ColumnWidth = 13.4101848602295
